$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter group member names (column A) and emails (column B).
# Order chosen so shared-string table is built up in the same sequence
# as the target workbook.
$ws.Range("A4").Value = "Mateusz Kuzon"
$ws.Range("B4").Value = "c3076381@student.shu.ac.uk"
$ws.Range("A5").Value = "Scott Forster"
$ws.Range("A6").Value = "Charlie Molnar"
$ws.Range("A7").Value = "Nasir Al-Monsur"
$ws.Range("B7").Value = "c3077195@student.shu.ac.uk"
$ws.Range("B5").Value = "c3022458@student.shu.ac.uk"
$ws.Range("B6").Value = "c3035956@student.shu.ac.uk"

# Excel auto-detects the e-mail addresses and turns them into mailto
# hyperlinks (this also introduces the Hyperlink cell style / font).
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:c3076381@student.shu.ac.uk")
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:c3077195@student.shu.ac.uk")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:c3022458@student.shu.ac.uk")
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:c3035956@student.shu.ac.uk")

# Column B widened (it no longer shares column A's width) to fit the
# e-mail addresses.
$ws.Columns("B").ColumnWidth = 24.6

# Selection moved as last user action.
$null = $ws.Range("C16").Select()
